$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 corresponds to @cognizant_exam_discussion_groups
# D9: current_phase 1 -> 2
$ws.Range("D9").Value = 2

# E9: last_action_date "" -> timestamp
$ws.Range("E9").Value = "2026-02-19T11:42:18.853835+00:00"

# I9: replies_count 0 -> 1
$ws.Range("I9").Value = 1

# M9: replied_message_ids "[]" -> "[1960685]"
$ws.Range("M9").Value = "[1960685]"
